$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "alem"
$ws.Range("C3").Value = "tsehey"
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 911140743
$ws.Range("F3").Value = "GO"
$ws.Range("G3").Value = "ethiopia"

$ws.Range("A4").Select()
